# Automatische test-sync: 2025-06-24 22:07:50
# Adds the new incoming-mail log row (row 42) to the "Logs" sheet,
# expands the conditional-formatting ranges that covered the log table,
# and bumps the "Overig" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$wsLogs = $wb.Worksheets.Item("Logs")

# New row of data describing the incoming mail.
$wsLogs.Range("A42").Value = "Beschadigd product ontvangen"
$wsLogs.Range("B42").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C42").Value = "Het product dat ik heb ontvangen is beschadigd aangekomen."
$wsLogs.Range("D42").Value = "Overig"
$wsLogs.Range("F42").Value = "2025-06-24 22:06:53"
$wsLogs.Range("G42").Value = "Nee"

# Extend the existing conditional-formatting rules (Categorie/Beantwoord
# columns) so they keep covering the whole table, now through row 42.
$fcsCategorie = $wsLogs.Range("D2:D41").FormatConditions
for ($i = 1; $i -le $fcsCategorie.Count; $i++) {
    $fcsCategorie.Item($i).ModifyAppliesToRange($wsLogs.Range("D2:D42"))
}

$fcsBeantwoord = $wsLogs.Range("G2:G41").FormatConditions
for ($i = 1; $i -le $fcsBeantwoord.Count; $i++) {
    $fcsBeantwoord.Item($i).ModifyAppliesToRange($wsLogs.Range("G2:G42"))
}

# Update the Dashboard summary count for the "Overig" category.
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Range("B9").Value = 2
